# Update crypto price/volume figures for the "Sat Jan 21 19:33:36 UTC 2023" symbol-list refresh.
# Source data is stored as literal text (t="inlineStr") in the workbook, so every updated
# cell is forced to Text format before the new value is written; this preserves exact
# formatting (leading/trailing zeros, "%" suffix, very small decimals, etc.) instead of
# letting Excel auto-convert the strings into numbers/percentages.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new literal text value, taken from the diff.
$updates = [ordered]@{
    "D2" = "305.06"
    "E2" = "4.41%"
    "D3" = "35.69"
    "E3" = "14.38%"
    "D4" = "5.103"
    "E4" = "3.20%"
    "D5" = "0.07840"
    "E5" = "4.61%"
    "D6" = "2.250"
    "E6" = "-1.15%"
    "D7" = "8.108"
    "E7" = "3.94%"
    "D8" = "4.007"
    "E8" = "6.21%"
    "D9" = "0.9261"
    "E9" = "0.50%"
    "D10" = "0.09887"
    "E10" = "7.26%"
    "D11" = "0.1821"
    "D12" = "0.08727"
    "E12" = "4.50%"
    "D13" = "0.03421"
    "E13" = "4.41%"
    "D14" = "0.09948"
    "E14" = "0.19%"
    "D15" = "0.001485"
    "E15" = "-0.65%"
    "D16" = "0.005738"
    "E16" = "-0.19%"
    "D17" = "3.482"
    "E17" = "0.07%"
    "E18" = "-1.79%"
    "D19" = "0.3433"
    "E19" = "2.55%"
    "D20" = "0.1322"
    "E20" = "0.62%"
    "D21" = "4.556"
    "E21" = "11.25%"
    "D22" = "0.2238"
    "E22" = "6.62%"
    "D23" = "0.04685"
    "E23" = "3.19%"
    "E24" = "1.73%"
    "D25" = "0.004503"
    "E25" = "4.52%"
    "D26" = "0.0001299"
    "E26" = "-0.12%"
    "D27" = "0.0002703"
    "E27" = "-20.24%"
    "D39" = "0.01760"
    "E39" = "8.04%"
    "D40" = "0.04719"
    "E40" = "3.06%"
    "D41" = "0.008027"
    "E41" = "7.57%"
    "D42" = "0.1423"
    "D43" = "0.008563"
    "E43" = "-12.87%"
    "D44" = "0.002211"
    "E44" = "2.35%"
    "D45" = "0.009122"
    "E45" = "-6.76%"
    "D46" = "0.00006213"
    "E46" = "1.94%"
    "D47" = "0.00000000751"
    "E47" = "0.08%"
    "D48" = "5.673"
    "E48" = "122.38%"
    "D49" = "0.002693"
    "E49" = "34.92%"
    "D50" = "0.00002103"
    "E50" = "0.08%"
    "D51" = "0.0002003"
    "E51" = "0.08%"
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"   # Text format -> keep the exact string, no auto number/percent conversion
    $cell.Value = $updates[$ref]
}

